$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)
$shape.Table.ApplyStyle("{EE7C059D-57EC-4A08-AE5D-AD4FF56BC18C}")
